$d = $word.ActiveDocument

# 1. Replace "module CMT-303" with "one of the modules" in the first
#    body paragraph (the sentence about the CMT-303 team-working module).
$d.Content.Find.Execute("module CMT-303", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "one of the modules", 2)

# 2. Move the "_GoBack" bookmark from the end of the second paragraph to
#    the midpoint of "communication" (between "communicatio" and "n") in
#    the sentence we just edited, matching the author's last-edit caret
#    position.
$d.Bookmarks("_GoBack").Delete()

$found = $d.Content
$found.Find.Execute("communication")
$splitPos = $found.Start + 12
$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bm)
